$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J10").Value = 5
$ws.Range("M10").Value = "nan MPH"
$ws.Range("M12").Value = "nan°"
$ws.Range("J14").Value = "Herbst"
$ws.Range("M14").Value = "Undefined"
$ws.Range("J15").Value = "Right"
$ws.Range("M15").Value = "Undefined"
$ws.Range("J16").Value = "83-85 MPH"
$ws.Range("J17").Value = "SL,CB,FB,CH"
$ws.Range("J19").Value = 3
$ws.Range("M19").Value = "nan MPH"
$ws.Range("M21").Value = "nan°"
$ws.Range("M23").Value = "Undefined"
$ws.Range("M24").Value = "Undefined"
$ws.Range("J25").Value = "88-90 MPH"
$ws.Range("J26").Value = "CB,FB,CH"
$ws.Range("J28").Value = 7
$ws.Range("M28").Value = "nan MPH"
$ws.Range("J29").Value = 2
$ws.Range("M30").Value = "nan°"
$ws.Range("J32").Value = "Plum"
$ws.Range("M32").Value = "Undefined"
$ws.Range("M33").Value = "Undefined"
$ws.Range("J34").Value = "84-86 MPH"
$ws.Range("J35").Value = "SL,FB,CH"
$ws.Range("J37").Value = 9
$ws.Range("M37").Value = "77.19 MPH"
$ws.Range("M39").Value = "-16.33°"
$ws.Range("J41").Value = "Thompson"
$ws.Range("M41").Value = "Ground Ball"
$ws.Range("J42").Value = "Left"
$ws.Range("M42").Value = "Single"
$ws.Range("J43").Value = "84-84 MPH"
$ws.Range("J44").Value = "SL,FB,CH"
$ws.Range("J46").Value = 4
$ws.Range("M46").Value = "92.3 MPH"
$ws.Range("M48").Value = "-3.48°"
$ws.Range("M50").Value = "Ground Ball"
$ws.Range("M51").Value = "Single"
$ws.Range("J52").Value = "88-90 MPH"
$ws.Range("J53").Value = "CB,FB,CH"
$ws.Range("J61").Value = 6
$ws.Range("M61").Value = "62.56 MPH"
$ws.Range("J62").Value = 1
$ws.Range("M63").Value = "41.13°"
$ws.Range("J65").Value = "Herbst"
$ws.Range("M65").Value = "Line Drive"
$ws.Range("M66").Value = "Out"
$ws.Range("J67").Value = "83-85 MPH"
$ws.Range("J68").Value = "SL,CB,FB,CH"

Write-Output "done"
